$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.047.65"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "3.174.92"
$ws.Range("E3").Value = "  -4.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "590.78"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "134.54"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -5.90%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.172.80"
$ws.Range("E8").Value = "  -4.57%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.515"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -5.87%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "5.24"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("E12").Value = "  -3.29%  "
$ws.Range("E13").Value = "  -4.58%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.05"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "3.696.95"
$ws.Range("E15").Value = "  -4.49%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "3.176.73"
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").Value = "63.018.12"
$ws.Range("E18").Value = "  -1.38%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.12%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "461.94"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -3.94%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.89"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("E22").Value = "  -5.26%  "
$ws.Range("E23").Value = "  -4.10%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.47"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.83%  "
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -4.00%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.86"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -4.84%  "
$ws.Range("E30").Value = "  -6.45%  "
$ws.Range("E31").Value = "  -5.51%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "27.18"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -6.20%  "
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -6.80%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.03"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -6.49%  "
$ws.Range("E36").Value = "  -4.05%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "51.36"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -6.30%  "
$ws.Range("E39").Value = "  -2.86%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "402.85"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -6.66%  "
$ws.Range("E41").Value = "  -2.79%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.112"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.60"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -5.47%  "
$ws.Range("D44").Value = "2.801.95"
$ws.Range("E44").Value = "  -10.18%  "
$ws.Range("E45").Value = "  -5.79%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  -5.34%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "124.96"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  -3.97%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "34.35"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -5.66%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
